# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the
# bfe1da4d-8542-4909-8738-8be88bbbf1cb entry (row 4) on both the
# zh-cn and de-de localization-status sheets, reflecting a fresh handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-04 16:56:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-04 16:56:28"
